# Applies the per-class Leve profit recompute from the scheduled runner.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3023.2666
$ws.Range("I64").Value = 2870.4
$ws.Range("J64").Value = 3099.7
$ws.Range("K64").Value = 2870.4
$ws.Range("L64").Value = 3099.7
$ws.Range("M64").Value = -2622.4
$ws.Range("N64").Value = -3595.7
$ws.Range("H67").Value = 3023.2666
$ws.Range("I67").Value = 2870.4
$ws.Range("J67").Value = 3099.7
$ws.Range("K67").Value = 2870.4
$ws.Range("L67").Value = 3099.7
$ws.Range("M67").Value = -2012.4
$ws.Range("N67").Value = -4815.7
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H74").Value = 4016.8667
$ws.Range("I74").Value = 3698.6
$ws.Range("J74").Value = 4176
$ws.Range("K74").Value = 3698.6
$ws.Range("L74").Value = 4176
$ws.Range("M74").Value = -2762.6
$ws.Range("N74").Value = -6048
$ws.Range("H77").Value = 4016.8667
$ws.Range("I77").Value = 3698.6
$ws.Range("J77").Value = 4176
$ws.Range("K77").Value = 18493
$ws.Range("L77").Value = 20880
$ws.Range("M77").Value = -13813
$ws.Range("N77").Value = -30240
$ws.Range("H88").Value = 6012
$ws.Range("I88").Value = 6477.143
$ws.Range("J88").Value = 5605
$ws.Range("K88").Value = 6477.143
$ws.Range("L88").Value = 5605
$ws.Range("M88").Value = -6071.143
$ws.Range("N88").Value = -6417
$ws.Range("H91").Value = 6012
$ws.Range("I91").Value = 6477.143
$ws.Range("J91").Value = 5605
$ws.Range("K91").Value = 6477.143
$ws.Range("L91").Value = 5605
$ws.Range("M91").Value = -5073.143
$ws.Range("N91").Value = -8413
$ws.Range("H133").Value = 28031.428
$ws.Range("J133").Value = 28031.428
$ws.Range("L133").Value = 28031.428
$ws.Range("N133").Value = -38151.428
$ws.Range("H138").Value = 3589.9512
$ws.Range("I138").Value = 1803.7084
$ws.Range("J138").Value = 6111.706
$ws.Range("K138").Value = 5411.1252
$ws.Range("L138").Value = 18335.118
$ws.Range("M138").Value = -271.1252000000004
$ws.Range("N138").Value = -28615.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2358.9443
$ws.Range("I61").Value = 1676.1428
$ws.Range("J61").Value = 4748.75
$ws.Range("K61").Value = 1676.1428
$ws.Range("L61").Value = 4748.75
$ws.Range("M61").Value = -1464.1428
$ws.Range("N61").Value = -5172.75
$ws.Range("H136").Value = 2358.9443
$ws.Range("I136").Value = 1676.1428
$ws.Range("J136").Value = 4748.75
$ws.Range("K136").Value = 5028.428400000001
$ws.Range("L136").Value = 14246.25
$ws.Range("M136").Value = -2478.428400000001
$ws.Range("N136").Value = -19346.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 122
$ws.Range("I5").Value = 122
$ws.Range("K5").Value = 122
$ws.Range("M5").Value = -9
$ws.Range("H6").Value = 28666.666
$ws.Range("J6").Value = 28666.666
$ws.Range("L6").Value = 28666.666
$ws.Range("N6").Value = -28892.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 80092
$ws.Range("J54").Value = 80092
$ws.Range("L54").Value = 80092
$ws.Range("N54").Value = -81408
$ws.Range("H123").Value = 32000
$ws.Range("J123").Value = 32000
$ws.Range("L123").Value = 32000
$ws.Range("N123").Value = -41800
$ws.Range("H129").Value = 40856.715
$ws.Range("I129").Value = 29333.334
$ws.Range("J129").Value = 49499.25
$ws.Range("K129").Value = 29333.334
$ws.Range("L129").Value = 49499.25
$ws.Range("M129").Value = -24333.334
$ws.Range("N129").Value = -59499.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1457.3334
$ws.Range("I109").Value = 530
$ws.Range("J109").Value = 1766.4445
$ws.Range("K109").Value = 1590
$ws.Range("L109").Value = 5299.333500000001
$ws.Range("M109").Value = -550
$ws.Range("N109").Value = -7379.333500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 70030
$ws.Range("J48").Value = 70030
$ws.Range("L48").Value = 70030
$ws.Range("N48").Value = -71000
$ws.Range("H49").Value = 38584.285
$ws.Range("J49").Value = 38584.285
$ws.Range("L49").Value = 38584.285
$ws.Range("N49").Value = -38952.285
$ws.Range("H80").Value = 3155.3572
$ws.Range("I80").Value = 3297.4443
$ws.Range("J80").Value = 2899.6
$ws.Range("K80").Value = 3297.4443
$ws.Range("L80").Value = 2899.6
$ws.Range("M80").Value = -2299.4443
$ws.Range("N80").Value = -4895.6
$ws.Range("H83").Value = 3155.3572
$ws.Range("I83").Value = 3297.4443
$ws.Range("J83").Value = 2899.6
$ws.Range("K83").Value = 16487.2215
$ws.Range("L83").Value = 14498
$ws.Range("M83").Value = -11495.2215
$ws.Range("N83").Value = -24482
$ws.Range("H97").Value = 3170
$ws.Range("I97").Value = 1804
$ws.Range("J97").Value = 10000
$ws.Range("K97").Value = 1804
$ws.Range("L97").Value = 10000
$ws.Range("M97").Value = -1308
$ws.Range("N97").Value = -10992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 169984.33
$ws.Range("I23").Value = 335635.34
$ws.Range("J23").Value = 4333.3335
$ws.Range("K23").Value = 335635.34
$ws.Range("L23").Value = 4333.3335
$ws.Range("M23").Value = -335405.34
$ws.Range("N23").Value = -4793.3335
$ws.Range("H24").Value = 20605.6
$ws.Range("J24").Value = 25007
$ws.Range("L24").Value = 25007
$ws.Range("N24").Value = -25693
$ws.Range("H30").Value = 2000
$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 2000
$ws.Range("N30").Value = -2216
$ws.Range("H82").Value = 3029.077
$ws.Range("I82").Value = 2297.25
$ws.Range("J82").Value = 4200
$ws.Range("K82").Value = 2297.25
$ws.Range("L82").Value = 4200
$ws.Range("M82").Value = -1936.25
$ws.Range("N82").Value = -4922
$ws.Range("H85").Value = 3029.077
$ws.Range("I85").Value = 2297.25
$ws.Range("J85").Value = 4200
$ws.Range("K85").Value = 2297.25
$ws.Range("L85").Value = 4200
$ws.Range("M85").Value = -1049.25
$ws.Range("N85").Value = -6696
$ws.Range("H93").Value = 2570.923
$ws.Range("I93").Value = 1142.2
$ws.Range("J93").Value = 7333.3335
$ws.Range("K93").Value = 1142.2
$ws.Range("L93").Value = 7333.3335
$ws.Range("M93").Value = 105.8
$ws.Range("N93").Value = -9829.333500000001
$ws.Range("H100").Value = 2389.4443
$ws.Range("I100").Value = 1390.8334
$ws.Range("J100").Value = 4386.6665
$ws.Range("K100").Value = 1390.8334
$ws.Range("L100").Value = 4386.6665
$ws.Range("M100").Value = -849.8334
$ws.Range("N100").Value = -5468.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2433.3333
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2650
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2650
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -3898
$ws.Range("H65").Value = 2433.3333
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2650
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 13250
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -19490

